$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: write a value as literal TEXT (not auto-converted to a number /
# date / percent by Excel's smart-entry parsing). We stage the text in a
# scratch cell far outside the used range, force it to Text format, then
# paste-special just the VALUE into the destination so the destination's
# own number format / style is left completely untouched.
# ---------------------------------------------------------------------------
$helper = $ws.Range("Z1")
$helper.NumberFormat = "@"

function Set-TextValue([string]$targetAddr, [string]$text) {
    $helper.Value = $text
    $helper.Copy()
    $ws.Range($targetAddr).PasteSpecial(-4163)   # xlPasteValues
    $excel.CutCopyMode = 0
}

# ---------------------------------------------------------------------------
# Class Statistics block (K/L columns near top of sheet)
# ---------------------------------------------------------------------------
$ws.Range("L6").Value = 153
$ws.Range("L7").Value = 3
Set-TextValue "L9"  "48.1%"
Set-TextValue "L10" "73.3%"

# ---------------------------------------------------------------------------
# "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System" swaps
# (Recorded-By order flip on the already-recorded session rows)
# ---------------------------------------------------------------------------
$swapRows = 8,9,34,35,60,61,86,87,112,113,138,139,167,194,221,248,275,302
foreach ($r in $swapRows) {
    $ws.Range("G$r").Value = "dnasr281@gmail.com, System"
}

# ---------------------------------------------------------------------------
# Session rows that flip from "Not Recorded" (pink) to "Recorded" (green).
# Each needs: fill/font format copied from the prior (already-Recorded) row,
# plus G/H/I values filled in.
# ---------------------------------------------------------------------------
$recordedFlips = @(
    @{ Row = 15;  Src = 14;  G = "dnasr281@gmail.com"; H = "22/26"; I = "Recorded" },
    @{ Row = 41;  Src = 40;  G = "dnasr281@gmail.com"; H = "24/27"; I = "Recorded" },
    @{ Row = 67;  Src = 66;  G = "dnasr281@gmail.com"; H = "19/26"; I = "Recorded" },
    @{ Row = 93;  Src = 92;  G = "dnasr281@gmail.com"; H = "21/27"; I = "Recorded" },
    @{ Row = 119; Src = 118; G = "dnasr281@gmail.com"; H = "29/30"; I = "Recorded" },
    @{ Row = 145; Src = 144; G = "dnasr281@gmail.com"; H = "18/23"; I = "Recorded" }
)

foreach ($flip in $recordedFlips) {
    $srcRange = $ws.Range("A$($flip.Src):I$($flip.Src)")
    $dstRange = $ws.Range("A$($flip.Row):I$($flip.Row)")
    $srcRange.Copy()
    $dstRange.PasteSpecial(-4122)   # xlPasteFormats
    $excel.CutCopyMode = 0

    $ws.Range("G$($flip.Row)").Value = $flip.G
    $ws.Range("H$($flip.Row)").Value = $flip.H
    $ws.Range("I$($flip.Row)").Value = $flip.I
}

# ---------------------------------------------------------------------------
# Group Statistics block (rows 15-20, columns O/P/R/S) recalculated
# after the B1A1 session-15 recording above.
# ---------------------------------------------------------------------------
$ws.Range("O15").Value = 13
$ws.Range("P15").Value = 1
Set-TextValue "R15" "50.0%"
Set-TextValue "S15" "81.1%"

$ws.Range("O16").Value = 14
$ws.Range("P16").Value = 0
Set-TextValue "R16" "53.8%"
Set-TextValue "S16" "79.4%"

$ws.Range("O17").Value = 14
$ws.Range("P17").Value = 0
Set-TextValue "R17" "53.8%"
Set-TextValue "S17" "65.7%"

$ws.Range("O18").Value = 14
$ws.Range("P18").Value = 0
Set-TextValue "R18" "53.8%"
Set-TextValue "S18" "69.8%"

$ws.Range("O19").Value = 14
$ws.Range("P19").Value = 0
Set-TextValue "R19" "53.8%"
Set-TextValue "S19" "74.0%"

$ws.Range("O20").Value = 13
$ws.Range("P20").Value = 1
Set-TextValue "R20" "50.0%"
Set-TextValue "S20" "75.6%"

# ---------------------------------------------------------------------------
# Clean up the scratch helper cell so it leaves no trace.
# ---------------------------------------------------------------------------
$helper.Clear()
